$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 79, pushing current rows 79:90 down to 80:91.
$ws.Rows(79).Insert()

# Populate the new row 79 with the new data record.
$ws.Range("A79").Value2 = 4
$ws.Range("B79").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value2 = "Los Lagos"
$ws.Range("D79").Value2 = 44491
$ws.Range("E79").Value2 = 10
$ws.Range("F79").Value2 = "Fruta"
$ws.Range("G79").Value2 = 100108
$ws.Range("H79").Value2 = "Tropicales y subtropicales"
$ws.Range("I79").Value2 = 100108002
$ws.Range("J79").Value2 = "Mango"
$ws.Range("K79").Value2 = "Sin especificar"
$ws.Range("L79").Value2 = "Primera"
$ws.Range("M79").Value2 = 200
$ws.Range("N79").Value2 = 7000
$ws.Range("O79").Value2 = 7500
$ws.Range("P79").Value2 = 7250
$ws.Range("Q79").Value2 = "`$/bandeja 4 kilos"
$ws.Range("R79").Value2 = "Perú"
$ws.Range("S79").Value2 = 1812
$ws.Range("T79").Value2 = 4
